$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn (E) and de-de (F) status columns, rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:E4").Value = "In Translation"
$wsOverview.Range("F2:F4").Value = "In Translation"

# zh-cn sheet: Status column (C), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column (C), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Narrow the Status columns to fit the shorter text ---
# ColumnWidth 12.5 -> stored OOXML column width ~13.33 chars (closest
# achievable value to the narrower width produced by the text change)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
